# Update historique.xlsx: swap E/F header meaning, rewrite existing rows 2-8
# with new stock-movement data, and append new rows 9-16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: columns E and F swap meaning (Quantite_Avant / Quantite_Apres) ---
$ws.Cells.Item(1, 5).Value = "Quantite_Avant"
$ws.Cells.Item(1, 6).Value = "Quantite_Apres"

# --- Full data set for rows 2..16: Date, Produit, Nature, Quantite_Mouvement, Quantite_Avant, Quantite_Apres ---
$data = @(
    @("2025-05-23 06:40:52", "Tournevis cruciforme", "Sortie",     1,    103,  102),
    @("2025-05-23 07:10:26", "Clé à molette",        "Sortie",     2,    25,   23),
    @("2025-05-23 07:34:03", "Perceuse sans fil",    "Sortie",     10,   19,   9),
    @("2025-05-23 07:34:21", "Perceuse sans fil",    "Entrée",     2,    9,    11),
    @("2025-05-23 07:34:51", "Perceuse sans fil",    "Inventaire", 2,    11,   9),
    @("2025-05-23 07:36:07", "Vis 6x50mm",           "Sortie",     500,  991,  491),
    @("2025-05-23 07:36:29", "Vis 6x50mm",           "Entrée",     12,   491,  503),
    @("2025-05-23 07:37:13", "Vis 6x50mm",           "Entrée",     1000, 503,  1503),
    @("2025-05-23 07:37:45", "Vis 6x50mm",           "Sortie",     1000, 1503, 503),
    @("2025-05-23 07:40:59", "Clé à molette",        "Sortie",     12,   23,   11),
    @("2025-05-23 07:53:38", "Clé à molette",        "Entrée",     23,   11,   34),
    @("2025-05-23 07:53:53", "Clé à molette",        "Entrée",     8,    34,   42),
    @("2025-05-23 07:53:59", "Clé à molette",        "Entrée",     3,    42,   45),
    @("2025-05-23 07:54:04", "Clé à molette",        "Entrée",     1,    45,   46),
    @("2025-05-23 07:54:14", "Clé à molette",        "Sortie",     4,    46,   42)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
    $row = $row + 1
}
